# Regenerate merged AHB files
# - Rename header row labels from *_old/*_new to *_FV2310/*_FV2404
# - Turn the data range into an Excel Table (Table1)
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row (A1:U1) from *_old / *_new suffixes to the
#    FV2310 / FV2404 release names.
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Convert A1:U79 into a native Excel table ("Table1").
$tableRange = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row so row 1 stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
